$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3348.25
$ws.Range("J17").Value = 3367.6296
$ws.Range("L17").Value = 10102.8888
$ws.Range("N17").Value = -10438.8888

$ws.Range("H116").Value = 55565320
$ws.Range("J116").Value = 8872.75
$ws.Range("L116").Value = 8872.75
$ws.Range("N116").Value = -15756.75

$ws.Range("H132").Value = 1888.4286
$ws.Range("I132").Value = 1807.4667
$ws.Range("J132").Value = 2799.25
$ws.Range("K132").Value = 5422.4001
$ws.Range("L132").Value = 8397.75
$ws.Range("M132").Value = -2892.4001
$ws.Range("N132").Value = -13457.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38583596
$ws.Range("I32").Value = 52158390
$ws.Range("J32").Value = 6497716
$ws.Range("K32").Value = 52158390
$ws.Range("L32").Value = 6497716
$ws.Range("M32").Value = -52158103
$ws.Range("N32").Value = -6498290

$ws.Range("H45").Value = 2880.3333
$ws.Range("I45").Value = 1608.3334
$ws.Range("J45").Value = 4152.3335
$ws.Range("K45").Value = 1608.3334
$ws.Range("L45").Value = 4152.3335
$ws.Range("M45").Value = -1231.3334
$ws.Range("N45").Value = -4906.3335

$ws.Range("H61").Value = 2909.9033
$ws.Range("I61").Value = 2515.35
$ws.Range("K61").Value = 2515.35
$ws.Range("M61").Value = -2303.35

$ws.Range("H102").Value = 1467.7858
$ws.Range("I102").Value = 1292.2307
$ws.Range("K102").Value = 1292.2307
$ws.Range("M102").Value = 329.7692999999999

$ws.Range("H136").Value = 2909.9033
$ws.Range("I136").Value = 2515.35
$ws.Range("K136").Value = 7546.049999999999
$ws.Range("M136").Value = -4996.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28126.105
$ws.Range("I20").Value = 46809.637
$ws.Range("J20").Value = 2436.25
$ws.Range("K20").Value = 46809.637
$ws.Range("L20").Value = 2436.25
$ws.Range("M20").Value = -46562.637
$ws.Range("N20").Value = -2930.25

$ws.Range("H80").Value = 536.875
$ws.Range("J80").Value = 611.3
$ws.Range("L80").Value = 611.3
$ws.Range("N80").Value = -2607.3

$ws.Range("H83").Value = 536.875
$ws.Range("J83").Value = 611.3
$ws.Range("L83").Value = 3056.5
$ws.Range("N83").Value = -13040.5

$ws.Range("H94").Value = 637.35486
$ws.Range("I94").Value = 535
$ws.Range("J94").Value = 852.3
$ws.Range("K94").Value = 535
$ws.Range("L94").Value = 852.3
$ws.Range("M94").Value = -84
$ws.Range("N94").Value = -1754.3

$ws.Range("H134").Value = 5498032.5
$ws.Range("I134").Value = 6496356.5
$ws.Range("K134").Value = 19489069.5
$ws.Range("M134").Value = -19486534.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5370.5
$ws.Range("I31").Value = 2437
$ws.Range("J31").Value = 6426.56
$ws.Range("K31").Value = 2437
$ws.Range("L31").Value = 6426.56
$ws.Range("M31").Value = -2142
$ws.Range("N31").Value = -7016.56

$ws.Range("H34").Value = 5370.5
$ws.Range("I34").Value = 2437
$ws.Range("J34").Value = 6426.56
$ws.Range("K34").Value = 2437
$ws.Range("L34").Value = 6426.56
$ws.Range("M34").Value = -2235
$ws.Range("N34").Value = -6830.56

$ws.Range("H99").Value = 3412
$ws.Range("I99").Value = 3412
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3412
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1914
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 3851113.5
$ws.Range("I122").Value = 5885800
$ws.Range("K122").Value = 17657400
$ws.Range("M122").Value = -17654950

$ws.Range("H126").Value = 3412
$ws.Range("I126").Value = 3412
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10236
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7766
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3526.4285
$ws.Range("I132").Value = 3041.4482
$ws.Range("J132").Value = 4608.3076
$ws.Range("K132").Value = 9124.3446
$ws.Range("L132").Value = 13824.9228
$ws.Range("M132").Value = -6594.3446
$ws.Range("N132").Value = -18884.9228

$ws.Range("H134").Value = 3437.5
$ws.Range("I134").Value = 3416.6667
$ws.Range("K134").Value = 10250.0001
$ws.Range("M134").Value = -7715.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 305.2
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 281.5
$ws.Range("K2").Value = 2400
$ws.Range("L2").Value = 1689
$ws.Range("M2").Value = -2287
$ws.Range("N2").Value = -1915

$ws.Range("H5").Value = 2444.25
$ws.Range("I5").Value = 1558
$ws.Range("J5").Value = 3330.5
$ws.Range("K5").Value = 4674
$ws.Range("L5").Value = 9991.5
$ws.Range("M5").Value = -4562
$ws.Range("N5").Value = -10215.5

$ws.Range("H135").Value = 2444.25
$ws.Range("I135").Value = 1558
$ws.Range("J135").Value = 3330.5
$ws.Range("K135").Value = 14022
$ws.Range("L135").Value = 29974.5
$ws.Range("M135").Value = -11487
$ws.Range("N135").Value = -35044.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 11165.444
$ws.Range("I2").Value = 65.57143000000001
$ws.Range("J2").Value = 50015
$ws.Range("K2").Value = 65.57143000000001
$ws.Range("L2").Value = 50015
$ws.Range("M2").Value = 47.42856999999999
$ws.Range("N2").Value = -50241

$ws.Range("H70").Value = 26635.285
$ws.Range("I70").Value = 91398.28999999999
$ws.Range("K70").Value = 91398.28999999999
$ws.Range("M70").Value = -91128.28999999999

$ws.Range("H73").Value = 26635.285
$ws.Range("I73").Value = 91398.28999999999
$ws.Range("K73").Value = 91398.28999999999
$ws.Range("M73").Value = -90462.28999999999

$ws.Range("H80").Value = 3840
$ws.Range("I80").Value = 3815.6667
$ws.Range("K80").Value = 3815.6667
$ws.Range("M80").Value = -2817.6667

$ws.Range("H83").Value = 3840
$ws.Range("I83").Value = 3815.6667
$ws.Range("K83").Value = 19078.3335
$ws.Range("M83").Value = -14086.3335

$ws.Range("H97").Value = 639.3103599999999
$ws.Range("I97").Value = 417.08334
$ws.Range("K97").Value = 417.08334
$ws.Range("M97").Value = 78.91665999999998

$ws.Range("H102").Value = 2215
$ws.Range("I102").Value = 2555.9
$ws.Range("J102").Value = 1078.6666
$ws.Range("K102").Value = 2555.9
$ws.Range("L102").Value = 1078.6666
$ws.Range("M102").Value = -933.9000000000001
$ws.Range("N102").Value = -4322.6666

$ws.Range("H107").Value = 599.4545000000001
$ws.Range("I107").Value = 549.4286
$ws.Range("J107").Value = 687
$ws.Range("K107").Value = 549.4286
$ws.Range("L107").Value = 687
$ws.Range("M107").Value = 1370.5714
$ws.Range("N107").Value = -4527

$ws.Range("H113").Value = 12825.223
$ws.Range("J113").Value = 18170.166
$ws.Range("L113").Value = 18170.166
$ws.Range("N113").Value = -22510.166

$ws.Range("H134").Value = 151249.75
$ws.Range("J134").Value = 151249.75
$ws.Range("L134").Value = 453749.25
$ws.Range("N134").Value = -458819.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11889.667
$ws.Range("I7").Value = 11889.667
$ws.Range("K7").Value = 11889.667
$ws.Range("M7").Value = -11777.667

$ws.Range("H126").Value = 11889.667
$ws.Range("I126").Value = 11889.667
$ws.Range("K126").Value = 35669.001
$ws.Range("M126").Value = -33199.001

$ws.Range("H132").Value = 718416.5
$ws.Range("I132").Value = 1429941.2
$ws.Range("J132").Value = 6891.7144
$ws.Range("K132").Value = 4289823.6
$ws.Range("L132").Value = 20675.1432
$ws.Range("M132").Value = -4287293.6
$ws.Range("N132").Value = -25735.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 22683.637
$ws.Range("I4").Value = 24916
$ws.Range("K4").Value = 24916
$ws.Range("M4").Value = -24803

$ws.Range("H122").Value = 90916780
$ws.Range("I122").Value = 125008710
$ws.Range("K122").Value = 375026130
$ws.Range("M122").Value = -375023680

$ws.Range("H132").Value = 40959.04
$ws.Range("I132").Value = 44163.793
$ws.Range("J132").Value = 2502
$ws.Range("K132").Value = 132491.379
$ws.Range("L132").Value = 7506
$ws.Range("M132").Value = -129961.379
$ws.Range("N132").Value = -12566

$ws.Range("H136").Value = 2182.818
$ws.Range("I136").Value = 1471.9412
$ws.Range("J136").Value = 4599.8
$ws.Range("K136").Value = 4415.8236
$ws.Range("L136").Value = 13799.4
$ws.Range("M136").Value = -1865.8236
$ws.Range("N136").Value = -18899.4
